# The "invalid" sheet's B3 cell previously held the formula "=B2+1"
# (cached value 20241002). The edit converts it into a hard-coded
# literal value of 20241001. B4:B8 keep their shared formula
# ("=B3+1" chained via si="0"), so their cached results ripple forward
# automatically once B3's value changes (recalculated after the script
# runs): B4->20241002, B5->20241003, B6->20241004, B7->20241005,
# B8->20241006.
#
# The sheet view's selection also moves from K11 to B3 (and the
# scrolled-to topLeftCell is cleared as a result).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("invalid")

# Replace the formula in B3 with a plain literal value (removes the
# formula entirely, same as typing a bare number over it in Excel).
$ws.Range("B3").Value = 20241001

# Update the active selection/cell to B3.
$ws.Range("B3").Select()
